$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period label text
$ws.Range("B2").Value = "Periode : [Month description] [year]"

# Move selection to B2 (as saved in the file)
$ws.Range("B2").Select()
